# The template ships with a single placeholder slide (Title + Subtitle)
# that was only ever used as scaffolding while the deck's master/layouts
# were authored. Slides for this template are generated programmatically
# by the consuming app, so the shipped template itself should contain no
# slides. Remove the lone slide: this drops it from the slide id list and
# removes its backing part from the package.
$p = $ppt.ActivePresentation

while ($p.Slides.Count -gt 0) {
    $p.Slides.Item(1).Delete()
}
